$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1280
$ws.Range("I12").Value = 118.333336
$ws.Range("J12").Value = 3022.5
$ws.Range("K12").Value = 118.333336
$ws.Range("L12").Value = 3022.5
$ws.Range("M12").Value = 51.666664
$ws.Range("N12").Value = -3362.5
$ws.Range("H15").Value = 1226.2354
$ws.Range("I15").Value = 1226.2354
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3678.7062
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3509.7062
$ws.Range("H17").Value = 294556.5
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 299303.38
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 897910.14
$ws.Range("M17").Value = -582
$ws.Range("N17").Value = -898246.14
$ws.Range("H34").Value = 2666.3333
$ws.Range("I34").Value = 2666.3333
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2666.3333
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2463.3333
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 2666.3333
$ws.Range("I36").Value = 2666.3333
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2666.3333
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1951.3333
$ws.Range("N36").ClearContents()
$ws.Range("H38").Value = 59
$ws.Range("I38").Value = 50.727272
$ws.Range("J38").Value = 150
$ws.Range("K38").Value = 152.181816
$ws.Range("L38").Value = 450
$ws.Range("M38").Value = 219.818184
$ws.Range("N38").Value = -1194
$ws.Range("H80").Value = 2282.2942
$ws.Range("I80").Value = 2043
$ws.Range("J80").Value = 2449.8
$ws.Range("K80").Value = 6129
$ws.Range("L80").Value = 7349.400000000001
$ws.Range("M80").Value = -5131
$ws.Range("N80").Value = -9345.400000000001
$ws.Range("H83").Value = 2282.2942
$ws.Range("I83").Value = 2043
$ws.Range("J83").Value = 2449.8
$ws.Range("K83").Value = 18387
$ws.Range("L83").Value = 22048.2
$ws.Range("M83").Value = -13395
$ws.Range("N83").Value = -32032.2
$ws.Range("H86").Value = 8142.1816
$ws.Range("I86").Value = 7832.5
$ws.Range("J86").Value = 8211
$ws.Range("K86").Value = 7832.5
$ws.Range("L86").Value = 8211
$ws.Range("M86").Value = -6709.5
$ws.Range("N86").Value = -10457
$ws.Range("H89").Value = 8142.1816
$ws.Range("I89").Value = 7832.5
$ws.Range("J89").Value = 8211
$ws.Range("K89").Value = 39162.5
$ws.Range("L89").Value = 41055
$ws.Range("M89").Value = -33546.5
$ws.Range("N89").Value = -52287
$ws.Range("H96").Value = 606.1429000000001
$ws.Range("I96").Value = 606.1429000000001
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1818.4287
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -445.4287000000002
$ws.Range("H101").Value = 350
$ws.Range("I101").Value = 325
$ws.Range("J101").Value = 375
$ws.Range("K101").Value = 975
$ws.Range("L101").Value = 1125
$ws.Range("M101").Value = 647
$ws.Range("N101").Value = -4369
$ws.Range("H116").Value = 9799.478999999999
$ws.Range("I116").Value = 11189.637
$ws.Range("J116").Value = 8525.166999999999
$ws.Range("K116").Value = 11189.637
$ws.Range("L116").Value = 8525.166999999999
$ws.Range("M116").Value = -7747.637000000001
$ws.Range("N116").Value = -15409.167
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = 0
$ws.Range("H129").Value = 1782
$ws.Range("I129").Value = 2654.8
$ws.Range("J129").Value = 327.33334
$ws.Range("K129").Value = 7964.400000000001
$ws.Range("L129").Value = 982.0000200000001
$ws.Range("M129").Value = -2964.400000000001
$ws.Range("N129").Value = -10982.00002
$ws.Range("H132").Value = 17545620
$ws.Range("I132").Value = 19609582
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 58828746
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -58826216
$ws.Range("N132").Value = -10910
$ws.Range("H135").Value = 1151.3158
$ws.Range("I135").Value = 363.69232
$ws.Range("J135").Value = 2857.8333
$ws.Range("K135").Value = 3273.23088
$ws.Range("L135").Value = 25720.4997
$ws.Range("M135").Value = -738.2308800000001
$ws.Range("N135").Value = -30790.4997
$ws.Range("H138").Value = 3242.28
$ws.Range("I138").Value = 1318.5
$ws.Range("J138").Value = 5690.727
$ws.Range("K138").Value = 3955.5
$ws.Range("L138").Value = 17072.181
$ws.Range("M138").Value = 1184.5
$ws.Range("N138").Value = -27352.181
$ws.Range("H141").Value = 1506.92
$ws.Range("I141").Value = 1464.0435
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 4392.1305
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 787.8694999999998
$ws.Range("N141").Value = -16360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3384.6155
$ws.Range("I2").Value = 3414.0908
$ws.Range("J2").Value = 3363
$ws.Range("K2").Value = 3414.0908
$ws.Range("L2").Value = 3363
$ws.Range("M2").Value = -3301.0908
$ws.Range("N2").Value = -3589
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").Value = 0
$ws.Range("H32").Value = 4044.4883
$ws.Range("I32").Value = 3467.0256
$ws.Range("J32").Value = 9674.75
$ws.Range("K32").Value = 3467.0256
$ws.Range("L32").Value = 9674.75
$ws.Range("M32").Value = -3180.0256
$ws.Range("N32").Value = -10248.75
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H45").Value = 10844.077
$ws.Range("I45").Value = 19329.834
$ws.Range("J45").Value = 3570.5715
$ws.Range("K45").Value = 19329.834
$ws.Range("L45").Value = 3570.5715
$ws.Range("M45").Value = -18952.834
$ws.Range("N45").Value = -4324.5715
$ws.Range("H61").Value = 5769.5
$ws.Range("I61").Value = 5622.75
$ws.Range("J61").Value = 6063
$ws.Range("K61").Value = 5622.75
$ws.Range("L61").Value = 6063
$ws.Range("M61").Value = -5410.75
$ws.Range("N61").Value = -6487
$ws.Range("H62").Value = 60000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 60000
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 60000
$ws.Range("N62").Value = -61248
$ws.Range("H65").Value = 60000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 60000
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 180000
$ws.Range("N65").Value = -186240
$ws.Range("H92").Value = 55772.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 55772.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 55772.5
$ws.Range("N92").Value = -60764.5
$ws.Range("H110").Value = 9323.134
$ws.Range("I110").Value = 16133
$ws.Range("J110").Value = 4783.222
$ws.Range("K110").Value = 16133
$ws.Range("L110").Value = 4783.222
$ws.Range("M110").Value = -14088
$ws.Range("N110").Value = -8873.222
$ws.Range("H116").Value = 3384.6155
$ws.Range("I116").Value = 3414.0908
$ws.Range("J116").Value = 3363
$ws.Range("K116").Value = 3414.0908
$ws.Range("L116").Value = 3363
$ws.Range("M116").Value = -1120.0908
$ws.Range("N116").Value = -7951
$ws.Range("H132").Value = 4396
$ws.Range("I132").Value = 3799.8
$ws.Range("J132").Value = 5389.6665
$ws.Range("K132").Value = 11399.4
$ws.Range("L132").Value = 16168.9995
$ws.Range("M132").Value = -8869.400000000001
$ws.Range("N132").Value = -21228.9995
$ws.Range("H136").Value = 5769.5
$ws.Range("I136").Value = 5622.75
$ws.Range("J136").Value = 6063
$ws.Range("K136").Value = 16868.25
$ws.Range("L136").Value = 18189
$ws.Range("M136").Value = -14318.25
$ws.Range("N136").Value = -23289

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3384.6155
$ws.Range("I3").Value = 3414.0908
$ws.Range("J3").Value = 3363
$ws.Range("K3").Value = 3414.0908
$ws.Range("L3").Value = 3363
$ws.Range("M3").Value = -3300.0908
$ws.Range("N3").Value = -3591
$ws.Range("H12").Value = 385
$ws.Range("I12").Value = 170
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 170
$ws.Range("L12").Value = 600
$ws.Range("M12").Value = -2
$ws.Range("N12").Value = -936
$ws.Range("H20").Value = 4721.2666
$ws.Range("I20").Value = 6017.6
$ws.Range("J20").Value = 4073.1
$ws.Range("K20").Value = 6017.6
$ws.Range("L20").Value = 4073.1
$ws.Range("M20").Value = -5770.6
$ws.Range("N20").Value = -4567.1
$ws.Range("H94").Value = 1507.3334
$ws.Range("I94").Value = 1375.1034
$ws.Range("J94").Value = 2055.1428
$ws.Range("K94").Value = 1375.1034
$ws.Range("L94").Value = 2055.1428
$ws.Range("M94").Value = -924.1034
$ws.Range("N94").Value = -2957.1428
$ws.Range("H105").Value = 5481.6665
$ws.Range("I105").Value = 6797.5
$ws.Range("J105").Value = 2850
$ws.Range("K105").Value = 6797.5
$ws.Range("L105").Value = 2850
$ws.Range("M105").Value = -5050.5
$ws.Range("N105").Value = -6344
$ws.Range("H107").Value = 1325.1305
$ws.Range("I107").Value = 1153.6111
$ws.Range("J107").Value = 1942.6
$ws.Range("K107").Value = 1153.6111
$ws.Range("L107").Value = 1942.6
$ws.Range("M107").Value = 766.3888999999999
$ws.Range("N107").Value = -5782.6
$ws.Range("H134").Value = 2242.4412
$ws.Range("I134").Value = 2034.2413
$ws.Range("J134").Value = 3450
$ws.Range("K134").Value = 6102.7239
$ws.Range("L134").Value = 10350
$ws.Range("M134").Value = -3567.7239
$ws.Range("N134").Value = -15420

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1432.7778
$ws.Range("I22").Value = 240.83333
$ws.Range("J22").Value = 3816.6667
$ws.Range("K22").Value = 240.83333
$ws.Range("L22").Value = 3816.6667
$ws.Range("M22").Value = 109.16667
$ws.Range("N22").Value = -4516.6667
$ws.Range("H57").Value = 49166.668
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 49166.668
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 49166.668
$ws.Range("N57").Value = -50286.668
$ws.Range("H58").Value = 5125
$ws.Range("I58").Value = 4833.3335
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 4833.3335
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -4630.3335
$ws.Range("N58").Value = -6406
$ws.Range("H59").Value = 25727.273
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 25727.273
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 25727.273
$ws.Range("N59").Value = -28017.273
$ws.Range("H63").Value = 46250
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 46250
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 46250
$ws.Range("N63").Value = -47622
$ws.Range("H66").Value = 46250
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 46250
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 138750
$ws.Range("N66").Value = -145614
$ws.Range("H94").Value = 2024
$ws.Range("I94").Value = 2324.3333
$ws.Range("J94").Value = 222
$ws.Range("K94").Value = 2324.3333
$ws.Range("L94").Value = 222
$ws.Range("M94").Value = -1873.3333
$ws.Range("N94").Value = -1124
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H105").Value = 2179.8
$ws.Range("I105").Value = 2133.3333
$ws.Range("J105").Value = 2249.5
$ws.Range("K105").Value = 2133.3333
$ws.Range("L105").Value = 2249.5
$ws.Range("M105").Value = -386.3332999999998
$ws.Range("N105").Value = -5743.5
$ws.Range("H122").Value = 1750
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -9400
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("H134").Value = 45295.273
$ws.Range("I134").Value = 27281
$ws.Range("J134").Value = 93333.336
$ws.Range("K134").Value = 81843
$ws.Range("L134").Value = 280000.008
$ws.Range("M134").Value = -79308
$ws.Range("N134").Value = -285070.008
$ws.Range("H136").Value = 5125
$ws.Range("I136").Value = 4833.3335
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 14500.0005
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -11950.0005
$ws.Range("N136").Value = -23100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 272.51514
$ws.Range("I6").Value = 285.25806
$ws.Range("J6").Value = 75
$ws.Range("K6").Value = 855.77418
$ws.Range("L6").Value = 225
$ws.Range("M6").Value = -742.77418
$ws.Range("N6").Value = -451
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 30000
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -34868
$ws.Range("H128").Value = 350000
$ws.Range("I128").Value = 350000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 1050000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -1045020
$ws.Range("H134").Value = 5977.7144
$ws.Range("I134").Value = 4325.4116
$ws.Range("J134").Value = 13000
$ws.Range("K134").Value = 12976.2348
$ws.Range("L134").Value = 39000
$ws.Range("M134").Value = -7906.234800000002
$ws.Range("N134").Value = -49140
$ws.Range("H140").Value = 9920.333000000001
$ws.Range("I140").Value = 9920.333000000001
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 29760.999
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -24580.999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 283.9091
$ws.Range("I2").Value = 273.83334
$ws.Range("J2").Value = 310.77777
$ws.Range("K2").Value = 273.83334
$ws.Range("L2").Value = 310.77777
$ws.Range("M2").Value = -160.83334
$ws.Range("N2").Value = -536.7777699999999
$ws.Range("H38").Value = 12836.467
$ws.Range("I38").Value = 12272.728
$ws.Range("J38").Value = 14386.75
$ws.Range("K38").Value = 12272.728
$ws.Range("L38").Value = 14386.75
$ws.Range("M38").Value = -11809.728
$ws.Range("N38").Value = -15312.75
$ws.Range("H80").Value = 4767.75
$ws.Range("I80").Value = 2770.4
$ws.Range("J80").Value = 8096.6665
$ws.Range("K80").Value = 2770.4
$ws.Range("L80").Value = 8096.6665
$ws.Range("M80").Value = -1772.4
$ws.Range("N80").Value = -10092.6665
$ws.Range("H83").Value = 4767.75
$ws.Range("I83").Value = 2770.4
$ws.Range("J83").Value = 8096.6665
$ws.Range("K83").Value = 13852
$ws.Range("L83").Value = 40483.3325
$ws.Range("M83").Value = -8860
$ws.Range("N83").Value = -50467.3325
$ws.Range("H97").Value = 1018.93335
$ws.Range("I97").Value = 1018.93335
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1018.93335
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -522.93335
$ws.Range("H102").Value = 3105.2222
$ws.Range("I102").Value = 3118.375
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 3118.375
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1496.375
$ws.Range("N102").Value = -6244
$ws.Range("H104").Value = 30500
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 30500
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 30500
$ws.Range("N104").Value = -37488
$ws.Range("H132").Value = 2996.9167
$ws.Range("I132").Value = 2633
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 7899
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -5369
$ws.Range("N132").Value = -26060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2228.3333
$ws.Range("I46").Value = 1830
$ws.Range("J46").Value = 2626.6667
$ws.Range("K46").Value = 1830
$ws.Range("L46").Value = 2626.6667
$ws.Range("M46").Value = -1642
$ws.Range("N46").Value = -3002.6667
$ws.Range("H61").Value = 3079.3333
$ws.Range("I61").Value = 3464.4707
$ws.Range("J61").Value = 2144
$ws.Range("K61").Value = 3464.4707
$ws.Range("L61").Value = 2144
$ws.Range("M61").Value = -3262.4707
$ws.Range("N61").Value = -2548
$ws.Range("H93").Value = 2383.3809
$ws.Range("I93").Value = 2329
$ws.Range("J93").Value = 2900
$ws.Range("K93").Value = 2329
$ws.Range("L93").Value = 2900
$ws.Range("M93").Value = -1081
$ws.Range("N93").Value = -5396
$ws.Range("H113").Value = 3079.3333
$ws.Range("I113").Value = 3464.4707
$ws.Range("J113").Value = 2144
$ws.Range("K113").Value = 3464.4707
$ws.Range("L113").Value = 2144
$ws.Range("M113").Value = -1294.4707
$ws.Range("N113").Value = -6484
$ws.Range("H132").Value = 2616.5
$ws.Range("I132").Value = 2217.5312
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 6652.5936
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -4122.5936
$ws.Range("N132").Value = -32060
$ws.Range("H136").Value = 5500
$ws.Range("I136").Value = 5250
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 15750
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -13200
$ws.Range("N136").Value = -23700

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("H11").Value = 1000000
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -999858
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H62").Value = 14941.117
$ws.Range("I62").Value = 38000
$ws.Range("J62").Value = 9999.929
$ws.Range("K62").Value = 38000
$ws.Range("L62").Value = 9999.929
$ws.Range("M62").Value = -37376
$ws.Range("N62").Value = -11247.929
$ws.Range("H65").Value = 14941.117
$ws.Range("I65").Value = 38000
$ws.Range("J65").Value = 9999.929
$ws.Range("K65").Value = 190000
$ws.Range("L65").Value = 49999.645
$ws.Range("M65").Value = -186880
$ws.Range("N65").Value = -56239.645
$ws.Range("H81").Value = 10933.2
$ws.Range("I81").Value = 27675.25
$ws.Range("J81").Value = 4845.1816
$ws.Range("K81").Value = 55350.5
$ws.Range("L81").Value = 9690.3632
$ws.Range("M81").Value = -54289.5
$ws.Range("N81").Value = -11812.3632
$ws.Range("H84").Value = 10933.2
$ws.Range("I84").Value = 27675.25
$ws.Range("J84").Value = 4845.1816
$ws.Range("K84").Value = 276752.5
$ws.Range("L84").Value = 48451.816
$ws.Range("M84").Value = -271448.5
$ws.Range("N84").Value = -59059.816
$ws.Range("H96").Value = 1687.25
$ws.Range("I96").Value = 1599.8
$ws.Range("J96").Value = 1833
$ws.Range("K96").Value = 1599.8
$ws.Range("L96").Value = 1833
$ws.Range("M96").Value = -226.8
$ws.Range("N96").Value = -4579
$ws.Range("H110").Value = 179916.67
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 179916.67
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 179916.67
$ws.Range("N110").Value = -188096.67
$ws.Range("H113").Value = 741.7273
$ws.Range("I113").Value = 723.3077
$ws.Range("J113").Value = 768.3333
$ws.Range("K113").Value = 2169.9231
$ws.Range("L113").Value = 2304.9999
$ws.Range("M113").Value = 0.07690000000002328
$ws.Range("N113").Value = -6644.9999
$ws.Range("H122").Value = 3335.3
$ws.Range("I122").Value = 3039.2778
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 9117.8334
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -6667.8334
$ws.Range("N122").Value = -22898.5
$ws.Range("H124").Value = 41999.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 41999.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 41999.5
$ws.Range("N124").Value = -51819.5
$ws.Range("H126").Value = 2731
$ws.Range("I126").Value = 2243
$ws.Range("J126").Value = 4195
$ws.Range("K126").Value = 6729
$ws.Range("L126").Value = 12585
$ws.Range("M126").Value = -4259
$ws.Range("N126").Value = -17525
$ws.Range("H132").Value = 4725.577
$ws.Range("I132").Value = 4377
$ws.Range("J132").Value = 6189.6
$ws.Range("K132").Value = 13131
$ws.Range("L132").Value = 18568.8
$ws.Range("M132").Value = -10601
$ws.Range("N132").Value = -23628.8
$ws.Range("H136").Value = 2501.4736
$ws.Range("I136").Value = 2724.4167
$ws.Range("J136").Value = 2119.2856
$ws.Range("K136").Value = 8173.250100000001
$ws.Range("L136").Value = 6357.8568
$ws.Range("M136").Value = -5623.250100000001
$ws.Range("N136").Value = -11457.8568
